# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 126,18

$data[0,0] = 10
$data[0,1] = 'Vega Modelo de Temuco'
$data[0,2] = 'La Araucanía'
$data[0,3] = 44467
$data[0,4] = 9
$data[0,5] = 100112037
$data[0,6] = 'Cebollín'
$data[0,7] = 'Sin especificar'
$data[0,8] = 'Primera'
$data[0,9] = 50
$data[0,10] = 7000
$data[0,11] = 8000
$data[0,12] = 7600
$data[0,13] = '$/docena de paquetes'
$data[0,14] = 'Provincia de Cautín'
$data[0,15] = 633
$data[0,16] = 12
$data[0,17] = 'Hortaliza'

$data[1,0] = 10
$data[1,1] = 'Vega Modelo de Temuco'
$data[1,2] = 'La Araucanía'
$data[1,3] = 44467
$data[1,4] = 9
$data[1,5] = 100112037
$data[1,6] = 'Cebollín'
$data[1,7] = 'Sin especificar'
$data[1,8] = 'Primera'
$data[1,9] = 20
$data[1,10] = 5000
$data[1,11] = 5000
$data[1,12] = 5000
$data[1,13] = '$/docena de paquetes'
$data[1,14] = 'Región de O''Higgins'
$data[1,15] = 417
$data[1,16] = 12
$data[1,17] = 'Hortaliza'

$data[2,0] = 10
$data[2,1] = 'Vega Modelo de Temuco'
$data[2,2] = 'La Araucanía'
$data[2,3] = 44342
$data[2,4] = 9
$data[2,5] = 100112037
$data[2,6] = 'Cebollín'
$data[2,7] = 'Sin especificar'
$data[2,8] = 'Primera'
$data[2,9] = 50
$data[2,10] = 8000
$data[2,11] = 8000
$data[2,12] = 8000
$data[2,13] = '$/docena de paquetes'
$data[2,14] = 'Provincia de Cautín'
$data[2,15] = 667
$data[2,16] = 12
$data[2,17] = 'Hortaliza'

$data[3,0] = 10
$data[3,1] = 'Vega Modelo de Temuco'
$data[3,2] = 'La Araucanía'
$data[3,3] = 44434
$data[3,4] = 9
$data[3,5] = 100112037
$data[3,6] = 'Cebollín'
$data[3,7] = 'Sin especificar'
$data[3,8] = 'Primera'
$data[3,9] = 85
$data[3,10] = 5000
$data[3,11] = 5000
$data[3,12] = 5000
$data[3,13] = '$/docena de paquetes'
$data[3,14] = 'Región de O''Higgins'
$data[3,15] = 417
$data[3,16] = 12
$data[3,17] = 'Hortaliza'

$data[4,0] = 10
$data[4,1] = 'Vega Modelo de Temuco'
$data[4,2] = 'La Araucanía'
$data[4,3] = 44327
$data[4,4] = 9
$data[4,5] = 100112037
$data[4,6] = 'Cebollín'
$data[4,7] = 'Sin especificar'
$data[4,8] = 'Primera'
$data[4,9] = 70
$data[4,10] = 7000
$data[4,11] = 8000
$data[4,12] = 7714
$data[4,13] = '$/docena de paquetes'
$data[4,14] = 'Provincia de Cautín'
$data[4,15] = 643
$data[4,16] = 12
$data[4,17] = 'Hortaliza'

$data[5,0] = 10
$data[5,1] = 'Vega Modelo de Temuco'
$data[5,2] = 'La Araucanía'
$data[5,3] = 44327
$data[5,4] = 9
$data[5,5] = 100112037
$data[5,6] = 'Cebollín'
$data[5,7] = 'Sin especificar'
$data[5,8] = 'Primera'
$data[5,9] = 20
$data[5,10] = 8000
$data[5,11] = 8000
$data[5,12] = 8000
$data[5,13] = '$/docena de paquetes'
$data[5,14] = 'Región del Maule'
$data[5,15] = 667
$data[5,16] = 12
$data[5,17] = 'Hortaliza'

$data[6,0] = 10
$data[6,1] = 'Vega Modelo de Temuco'
$data[6,2] = 'La Araucanía'
$data[6,3] = 44300
$data[6,4] = 9
$data[6,5] = 100112037
$data[6,6] = 'Cebollín'
$data[6,7] = 'Sin especificar'
$data[6,8] = 'Primera'
$data[6,9] = 30
$data[6,10] = 4000
$data[6,11] = 4000
$data[6,12] = 4000
$data[6,13] = '$/docena de paquetes'
$data[6,14] = 'Provincia de Cautín'
$data[6,15] = 333
$data[6,16] = 12
$data[6,17] = 'Hortaliza'

$data[7,0] = 10
$data[7,1] = 'Vega Modelo de Temuco'
$data[7,2] = 'La Araucanía'
$data[7,3] = 44399
$data[7,4] = 9
$data[7,5] = 100112037
$data[7,6] = 'Cebollín'
$data[7,7] = 'Sin especificar'
$data[7,8] = 'Primera'
$data[7,9] = 100
$data[7,10] = 9000
$data[7,11] = 10000
$data[7,12] = 9500
$data[7,13] = '$/docena de paquetes'
$data[7,14] = 'Provincia de Cautín'
$data[7,15] = 792
$data[7,16] = 12
$data[7,17] = 'Hortaliza'

$data[8,0] = 10
$data[8,1] = 'Vega Modelo de Temuco'
$data[8,2] = 'La Araucanía'
$data[8,3] = 44399
$data[8,4] = 9
$data[8,5] = 100112037
$data[8,6] = 'Cebollín'
$data[8,7] = 'Sin especificar'
$data[8,8] = 'Primera'
$data[8,9] = 40
$data[8,10] = 5000
$data[8,11] = 5000
$data[8,12] = 5000
$data[8,13] = '$/docena de paquetes'
$data[8,14] = 'Región de O''Higgins'
$data[8,15] = 417
$data[8,16] = 12
$data[8,17] = 'Hortaliza'

$data[9,0] = 10
$data[9,1] = 'Vega Modelo de Temuco'
$data[9,2] = 'La Araucanía'
$data[9,3] = 44452
$data[9,4] = 9
$data[9,5] = 100112037
$data[9,6] = 'Cebollín'
$data[9,7] = 'Sin especificar'
$data[9,8] = 'Primera'
$data[9,9] = 40
$data[9,10] = 8000
$data[9,11] = 8000
$data[9,12] = 8000
$data[9,13] = '$/docena de paquetes'
$data[9,14] = 'Provincia de Cautín'
$data[9,15] = 667
$data[9,16] = 12
$data[9,17] = 'Hortaliza'

$data[10,0] = 10
$data[10,1] = 'Vega Modelo de Temuco'
$data[10,2] = 'La Araucanía'
$data[10,3] = 44452
$data[10,4] = 9
$data[10,5] = 100112037
$data[10,6] = 'Cebollín'
$data[10,7] = 'Sin especificar'
$data[10,8] = 'Primera'
$data[10,9] = 50
$data[10,10] = 5000
$data[10,11] = 5000
$data[10,12] = 5000
$data[10,13] = '$/docena de paquetes'
$data[10,14] = 'Región de O''Higgins'
$data[10,15] = 417
$data[10,16] = 12
$data[10,17] = 'Hortaliza'

$data[11,0] = 10
$data[11,1] = 'Vega Modelo de Temuco'
$data[11,2] = 'La Araucanía'
$data[11,3] = 44218
$data[11,4] = 9
$data[11,5] = 100112037
$data[11,6] = 'Cebollín'
$data[11,7] = 'Sin especificar'
$data[11,8] = 'Primera'
$data[11,9] = 40
$data[11,10] = 5000
$data[11,11] = 5000
$data[11,12] = 5000
$data[11,13] = '$/docena de paquetes'
$data[11,14] = 'Provincia de Cautín'
$data[11,15] = 417
$data[11,16] = 12
$data[11,17] = 'Hortaliza'

$data[12,0] = 10
$data[12,1] = 'Vega Modelo de Temuco'
$data[12,2] = 'La Araucanía'
$data[12,3] = 44383
$data[12,4] = 9
$data[12,5] = 100112037
$data[12,6] = 'Cebollín'
$data[12,7] = 'Sin especificar'
$data[12,8] = 'Primera'
$data[12,9] = 30
$data[12,10] = 8000
$data[12,11] = 8000
$data[12,12] = 8000
$data[12,13] = '$/docena de paquetes'
$data[12,14] = 'Provincia de Cautín'
$data[12,15] = 667
$data[12,16] = 12
$data[12,17] = 'Hortaliza'

$data[13,0] = 10
$data[13,1] = 'Vega Modelo de Temuco'
$data[13,2] = 'La Araucanía'
$data[13,3] = 44223
$data[13,4] = 9
$data[13,5] = 100112037
$data[13,6] = 'Cebollín'
$data[13,7] = 'Sin especificar'
$data[13,8] = 'Primera'
$data[13,9] = 65
$data[13,10] = 7000
$data[13,11] = 7000
$data[13,12] = 7000
$data[13,13] = '$/docena de paquetes'
$data[13,14] = 'Provincia de Cautín'
$data[13,15] = 583
$data[13,16] = 12
$data[13,17] = 'Hortaliza'

$data[14,0] = 10
$data[14,1] = 'Vega Modelo de Temuco'
$data[14,2] = 'La Araucanía'
$data[14,3] = 44441
$data[14,4] = 9
$data[14,5] = 100112037
$data[14,6] = 'Cebollín'
$data[14,7] = 'Sin especificar'
$data[14,8] = 'Primera'
$data[14,9] = 30
$data[14,10] = 9000
$data[14,11] = 9000
$data[14,12] = 9000
$data[14,13] = '$/docena de paquetes'
$data[14,14] = 'Provincia de Cautín'
$data[14,15] = 750
$data[14,16] = 12
$data[14,17] = 'Hortaliza'

$data[15,0] = 10
$data[15,1] = 'Vega Modelo de Temuco'
$data[15,2] = 'La Araucanía'
$data[15,3] = 44405
$data[15,4] = 9
$data[15,5] = 100112037
$data[15,6] = 'Cebollín'
$data[15,7] = 'Sin especificar'
$data[15,8] = 'Primera'
$data[15,9] = 55
$data[15,10] = 10000
$data[15,11] = 10000
$data[15,12] = 10000
$data[15,13] = '$/docena de paquetes'
$data[15,14] = 'Región de La Araucanía'
$data[15,15] = 833
$data[15,16] = 12
$data[15,17] = 'Hortaliza'

$data[16,0] = 10
$data[16,1] = 'Vega Modelo de Temuco'
$data[16,2] = 'La Araucanía'
$data[16,3] = 44405
$data[16,4] = 9
$data[16,5] = 100112037
$data[16,6] = 'Cebollín'
$data[16,7] = 'Sin especificar'
$data[16,8] = 'Primera'
$data[16,9] = 95
$data[16,10] = 5000
$data[16,11] = 5000
$data[16,12] = 5000
$data[16,13] = '$/docena de paquetes'
$data[16,14] = 'Región de O''Higgins'
$data[16,15] = 417
$data[16,16] = 12
$data[16,17] = 'Hortaliza'

$data[17,0] = 10
$data[17,1] = 'Vega Modelo de Temuco'
$data[17,2] = 'La Araucanía'
$data[17,3] = 44161
$data[17,4] = 9
$data[17,5] = 100112037
$data[17,6] = 'Cebollín'
$data[17,7] = 'Sin especificar'
$data[17,8] = 'Primera'
$data[17,9] = 50
$data[17,10] = 7000
$data[17,11] = 7000
$data[17,12] = 7000
$data[17,13] = '$/docena de paquetes'
$data[17,14] = 'Provincia de Cautín'
$data[17,15] = 583
$data[17,16] = 12
$data[17,17] = 'Hortaliza'

$data[18,0] = 10
$data[18,1] = 'Vega Modelo de Temuco'
$data[18,2] = 'La Araucanía'
$data[18,3] = 44161
$data[18,4] = 9
$data[18,5] = 100112037
$data[18,6] = 'Cebollín'
$data[18,7] = 'Sin especificar'
$data[18,8] = 'Primera'
$data[18,9] = 50
$data[18,10] = 5000
$data[18,11] = 5000
$data[18,12] = 5000
$data[18,13] = '$/docena de paquetes'
$data[18,14] = 'Región de O''Higgins'
$data[18,15] = 417
$data[18,16] = 12
$data[18,17] = 'Hortaliza'

$data[19,0] = 10
$data[19,1] = 'Vega Modelo de Temuco'
$data[19,2] = 'La Araucanía'
$data[19,3] = 44161
$data[19,4] = 9
$data[19,5] = 100112037
$data[19,6] = 'Cebollín'
$data[19,7] = 'Sin especificar'
$data[19,8] = 'Primera'
$data[19,9] = 40
$data[19,10] = 5000
$data[19,11] = 5000
$data[19,12] = 5000
$data[19,13] = '$/docena de paquetes'
$data[19,14] = 'Región del Maule'
$data[19,15] = 417
$data[19,16] = 12
$data[19,17] = 'Hortaliza'

$data[20,0] = 10
$data[20,1] = 'Vega Modelo de Temuco'
$data[20,2] = 'La Araucanía'
$data[20,3] = 44307
$data[20,4] = 9
$data[20,5] = 100112037
$data[20,6] = 'Cebollín'
$data[20,7] = 'Sin especificar'
$data[20,8] = 'Primera'
$data[20,9] = 45
$data[20,10] = 7000
$data[20,11] = 7000
$data[20,12] = 7000
$data[20,13] = '$/docena de paquetes'
$data[20,14] = 'Provincia de Cautín'
$data[20,15] = 583
$data[20,16] = 12
$data[20,17] = 'Hortaliza'

$data[21,0] = 10
$data[21,1] = 'Vega Modelo de Temuco'
$data[21,2] = 'La Araucanía'
$data[21,3] = 44193
$data[21,4] = 9
$data[21,5] = 100112037
$data[21,6] = 'Cebollín'
$data[21,7] = 'Sin especificar'
$data[21,8] = 'Primera'
$data[21,9] = 75
$data[21,10] = 6000
$data[21,11] = 6000
$data[21,12] = 6000
$data[21,13] = '$/docena de paquetes'
$data[21,14] = 'Provincia de Cautín'
$data[21,15] = 500
$data[21,16] = 12
$data[21,17] = 'Hortaliza'

$data[22,0] = 10
$data[22,1] = 'Vega Modelo de Temuco'
$data[22,2] = 'La Araucanía'
$data[22,3] = 44193
$data[22,4] = 9
$data[22,5] = 100112037
$data[22,6] = 'Cebollín'
$data[22,7] = 'Sin especificar'
$data[22,8] = 'Primera'
$data[22,9] = 55
$data[22,10] = 5000
$data[22,11] = 5000
$data[22,12] = 5000
$data[22,13] = '$/docena de paquetes'
$data[22,14] = 'Región de O''Higgins'
$data[22,15] = 417
$data[22,16] = 12
$data[22,17] = 'Hortaliza'

$data[23,0] = 10
$data[23,1] = 'Vega Modelo de Temuco'
$data[23,2] = 'La Araucanía'
$data[23,3] = 44336
$data[23,4] = 9
$data[23,5] = 100112037
$data[23,6] = 'Cebollín'
$data[23,7] = 'Sin especificar'
$data[23,8] = 'Primera'
$data[23,9] = 40
$data[23,10] = 7000
$data[23,11] = 8000
$data[23,12] = 7500
$data[23,13] = '$/docena de paquetes'
$data[23,14] = 'Provincia de Cautín'
$data[23,15] = 625
$data[23,16] = 12
$data[23,17] = 'Hortaliza'

$data[24,0] = 10
$data[24,1] = 'Vega Modelo de Temuco'
$data[24,2] = 'La Araucanía'
$data[24,3] = 44341
$data[24,4] = 9
$data[24,5] = 100112037
$data[24,6] = 'Cebollín'
$data[24,7] = 'Sin especificar'
$data[24,8] = 'Primera'
$data[24,9] = 40
$data[24,10] = 8000
$data[24,11] = 8000
$data[24,12] = 8000
$data[24,13] = '$/docena de paquetes'
$data[24,14] = 'Provincia de Cautín'
$data[24,15] = 667
$data[24,16] = 12
$data[24,17] = 'Hortaliza'

$data[25,0] = 10
$data[25,1] = 'Vega Modelo de Temuco'
$data[25,2] = 'La Araucanía'
$data[25,3] = 44426
$data[25,4] = 9
$data[25,5] = 100112037
$data[25,6] = 'Cebollín'
$data[25,7] = 'Sin especificar'
$data[25,8] = 'Primera'
$data[25,9] = 30
$data[25,10] = 7000
$data[25,11] = 7000
$data[25,12] = 7000
$data[25,13] = '$/docena de paquetes'
$data[25,14] = 'Provincia de Cautín'
$data[25,15] = 583
$data[25,16] = 12
$data[25,17] = 'Hortaliza'

$data[26,0] = 10
$data[26,1] = 'Vega Modelo de Temuco'
$data[26,2] = 'La Araucanía'
$data[26,3] = 44426
$data[26,4] = 9
$data[26,5] = 100112037
$data[26,6] = 'Cebollín'
$data[26,7] = 'Sin especificar'
$data[26,8] = 'Primera'
$data[26,9] = 40
$data[26,10] = 5000
$data[26,11] = 5000
$data[26,12] = 5000
$data[26,13] = '$/docena de paquetes'
$data[26,14] = 'Región de O''Higgins'
$data[26,15] = 417
$data[26,16] = 12
$data[26,17] = 'Hortaliza'

$data[27,0] = 10
$data[27,1] = 'Vega Modelo de Temuco'
$data[27,2] = 'La Araucanía'
$data[27,3] = 44328
$data[27,4] = 9
$data[27,5] = 100112037
$data[27,6] = 'Cebollín'
$data[27,7] = 'Sin especificar'
$data[27,8] = 'Primera'
$data[27,9] = 65
$data[27,10] = 8000
$data[27,11] = 8000
$data[27,12] = 8000
$data[27,13] = '$/docena de paquetes'
$data[27,14] = 'Provincia de Cautín'
$data[27,15] = 667
$data[27,16] = 12
$data[27,17] = 'Hortaliza'

$data[28,0] = 10
$data[28,1] = 'Vega Modelo de Temuco'
$data[28,2] = 'La Araucanía'
$data[28,3] = 44400
$data[28,4] = 9
$data[28,5] = 100112037
$data[28,6] = 'Cebollín'
$data[28,7] = 'Sin especificar'
$data[28,8] = 'Primera'
$data[28,9] = 30
$data[28,10] = 10000
$data[28,11] = 10000
$data[28,12] = 10000
$data[28,13] = '$/docena de paquetes'
$data[28,14] = 'Provincia de Cautín'
$data[28,15] = 833
$data[28,16] = 12
$data[28,17] = 'Hortaliza'

$data[29,0] = 10
$data[29,1] = 'Vega Modelo de Temuco'
$data[29,2] = 'La Araucanía'
$data[29,3] = 44392
$data[29,4] = 9
$data[29,5] = 100112037
$data[29,6] = 'Cebollín'
$data[29,7] = 'Sin especificar'
$data[29,8] = 'Primera'
$data[29,9] = 95
$data[29,10] = 5000
$data[29,11] = 5000
$data[29,12] = 5000
$data[29,13] = '$/docena de paquetes'
$data[29,14] = 'Región de O''Higgins'
$data[29,15] = 417
$data[29,16] = 12
$data[29,17] = 'Hortaliza'

$data[30,0] = 10
$data[30,1] = 'Vega Modelo de Temuco'
$data[30,2] = 'La Araucanía'
$data[30,3] = 44442
$data[30,4] = 9
$data[30,5] = 100112037
$data[30,6] = 'Cebollín'
$data[30,7] = 'Sin especificar'
$data[30,8] = 'Primera'
$data[30,9] = 30
$data[30,10] = 9000
$data[30,11] = 9000
$data[30,12] = 9000
$data[30,13] = '$/docena de paquetes'
$data[30,14] = 'Provincia de Cautín'
$data[30,15] = 750
$data[30,16] = 12
$data[30,17] = 'Hortaliza'

$data[31,0] = 10
$data[31,1] = 'Vega Modelo de Temuco'
$data[31,2] = 'La Araucanía'
$data[31,3] = 44453
$data[31,4] = 9
$data[31,5] = 100112037
$data[31,6] = 'Cebollín'
$data[31,7] = 'Sin especificar'
$data[31,8] = 'Primera'
$data[31,9] = 20
$data[31,10] = 8000
$data[31,11] = 8000
$data[31,12] = 8000
$data[31,13] = '$/docena de paquetes'
$data[31,14] = 'Provincia de Cautín'
$data[31,15] = 667
$data[31,16] = 12
$data[31,17] = 'Hortaliza'

$data[32,0] = 10
$data[32,1] = 'Vega Modelo de Temuco'
$data[32,2] = 'La Araucanía'
$data[32,3] = 44453
$data[32,4] = 9
$data[32,5] = 100112037
$data[32,6] = 'Cebollín'
$data[32,7] = 'Sin especificar'
$data[32,8] = 'Primera'
$data[32,9] = 20
$data[32,10] = 5000
$data[32,11] = 5000
$data[32,12] = 5000
$data[32,13] = '$/docena de paquetes'
$data[32,14] = 'Región de O''Higgins'
$data[32,15] = 417
$data[32,16] = 12
$data[32,17] = 'Hortaliza'

$data[33,0] = 10
$data[33,1] = 'Vega Modelo de Temuco'
$data[33,2] = 'La Araucanía'
$data[33,3] = 44309
$data[33,4] = 9
$data[33,5] = 100112037
$data[33,6] = 'Cebollín'
$data[33,7] = 'Sin especificar'
$data[33,8] = 'Primera'
$data[33,9] = 75
$data[33,10] = 6000
$data[33,11] = 7000
$data[33,12] = 6600
$data[33,13] = '$/docena de paquetes'
$data[33,14] = 'Provincia de Cautín'
$data[33,15] = 550
$data[33,16] = 12
$data[33,17] = 'Hortaliza'

$data[34,0] = 10
$data[34,1] = 'Vega Modelo de Temuco'
$data[34,2] = 'La Araucanía'
$data[34,3] = 44200
$data[34,4] = 9
$data[34,5] = 100112037
$data[34,6] = 'Cebollín'
$data[34,7] = 'Sin especificar'
$data[34,8] = 'Primera'
$data[34,9] = 30
$data[34,10] = 8000
$data[34,11] = 8000
$data[34,12] = 8000
$data[34,13] = '$/docena de paquetes'
$data[34,14] = 'Provincia de Cautín'
$data[34,15] = 667
$data[34,16] = 12
$data[34,17] = 'Hortaliza'

$data[35,0] = 10
$data[35,1] = 'Vega Modelo de Temuco'
$data[35,2] = 'La Araucanía'
$data[35,3] = 44208
$data[35,4] = 9
$data[35,5] = 100112037
$data[35,6] = 'Cebollín'
$data[35,7] = 'Sin especificar'
$data[35,8] = 'Primera'
$data[35,9] = 110
$data[35,10] = 8000
$data[35,11] = 8000
$data[35,12] = 8000
$data[35,13] = '$/docena de paquetes'
$data[35,14] = 'Provincia de Cautín'
$data[35,15] = 667
$data[35,16] = 12
$data[35,17] = 'Hortaliza'

$data[36,0] = 10
$data[36,1] = 'Vega Modelo de Temuco'
$data[36,2] = 'La Araucanía'
$data[36,3] = 44214
$data[36,4] = 9
$data[36,5] = 100112037
$data[36,6] = 'Cebollín'
$data[36,7] = 'Sin especificar'
$data[36,8] = 'Primera'
$data[36,9] = 50
$data[36,10] = 5000
$data[36,11] = 5000
$data[36,12] = 5000
$data[36,13] = '$/docena de paquetes'
$data[36,14] = 'Provincia de Cautín'
$data[36,15] = 417
$data[36,16] = 12
$data[36,17] = 'Hortaliza'

$data[37,0] = 10
$data[37,1] = 'Vega Modelo de Temuco'
$data[37,2] = 'La Araucanía'
$data[37,3] = 44354
$data[37,4] = 9
$data[37,5] = 100112037
$data[37,6] = 'Cebollín'
$data[37,7] = 'Sin especificar'
$data[37,8] = 'Primera'
$data[37,9] = 80
$data[37,10] = 9000
$data[37,11] = 9000
$data[37,12] = 9000
$data[37,13] = '$/docena de paquetes'
$data[37,14] = 'Provincia de Cautín'
$data[37,15] = 750
$data[37,16] = 12
$data[37,17] = 'Hortaliza'

$data[38,0] = 10
$data[38,1] = 'Vega Modelo de Temuco'
$data[38,2] = 'La Araucanía'
$data[38,3] = 44354
$data[38,4] = 9
$data[38,5] = 100112037
$data[38,6] = 'Cebollín'
$data[38,7] = 'Sin especificar'
$data[38,8] = 'Primera'
$data[38,9] = 50
$data[38,10] = 5000
$data[38,11] = 5000
$data[38,12] = 5000
$data[38,13] = '$/docena de paquetes'
$data[38,14] = 'Región Metropolitana'
$data[38,15] = 417
$data[38,16] = 12
$data[38,17] = 'Hortaliza'

$data[39,0] = 10
$data[39,1] = 'Vega Modelo de Temuco'
$data[39,2] = 'La Araucanía'
$data[39,3] = 44420
$data[39,4] = 9
$data[39,5] = 100112037
$data[39,6] = 'Cebollín'
$data[39,7] = 'Sin especificar'
$data[39,8] = 'Primera'
$data[39,9] = 155
$data[39,10] = 8000
$data[39,11] = 8000
$data[39,12] = 8000
$data[39,13] = '$/docena de paquetes'
$data[39,14] = 'Provincia de Cautín'
$data[39,15] = 667
$data[39,16] = 12
$data[39,17] = 'Hortaliza'

$data[40,0] = 10
$data[40,1] = 'Vega Modelo de Temuco'
$data[40,2] = 'La Araucanía'
$data[40,3] = 44420
$data[40,4] = 9
$data[40,5] = 100112037
$data[40,6] = 'Cebollín'
$data[40,7] = 'Sin especificar'
$data[40,8] = 'Primera'
$data[40,9] = 180
$data[40,10] = 5000
$data[40,11] = 5000
$data[40,12] = 5000
$data[40,13] = '$/docena de paquetes'
$data[40,14] = 'Región de O''Higgins'
$data[40,15] = 417
$data[40,16] = 12
$data[40,17] = 'Hortaliza'

$data[41,0] = 10
$data[41,1] = 'Vega Modelo de Temuco'
$data[41,2] = 'La Araucanía'
$data[41,3] = 44370
$data[41,4] = 9
$data[41,5] = 100112037
$data[41,6] = 'Cebollín'
$data[41,7] = 'Sin especificar'
$data[41,8] = 'Primera'
$data[41,9] = 20
$data[41,10] = 8000
$data[41,11] = 8000
$data[41,12] = 8000
$data[41,13] = '$/docena de paquetes'
$data[41,14] = 'Provincia de Cautín'
$data[41,15] = 667
$data[41,16] = 12
$data[41,17] = 'Hortaliza'

$data[42,0] = 10
$data[42,1] = 'Vega Modelo de Temuco'
$data[42,2] = 'La Araucanía'
$data[42,3] = 44370
$data[42,4] = 9
$data[42,5] = 100112037
$data[42,6] = 'Cebollín'
$data[42,7] = 'Sin especificar'
$data[42,8] = 'Primera'
$data[42,9] = 40
$data[42,10] = 5000
$data[42,11] = 5000
$data[42,12] = 5000
$data[42,13] = '$/docena de paquetes'
$data[42,14] = 'Región Metropolitana'
$data[42,15] = 417
$data[42,16] = 12
$data[42,17] = 'Hortaliza'

$data[43,0] = 10
$data[43,1] = 'Vega Modelo de Temuco'
$data[43,2] = 'La Araucanía'
$data[43,3] = 44237
$data[43,4] = 9
$data[43,5] = 100112037
$data[43,6] = 'Cebollín'
$data[43,7] = 'Sin especificar'
$data[43,8] = 'Primera'
$data[43,9] = 65
$data[43,10] = 7000
$data[43,11] = 7000
$data[43,12] = 7000
$data[43,13] = '$/docena de paquetes'
$data[43,14] = 'Provincia de Cautín'
$data[43,15] = 583
$data[43,16] = 12
$data[43,17] = 'Hortaliza'

$data[44,0] = 10
$data[44,1] = 'Vega Modelo de Temuco'
$data[44,2] = 'La Araucanía'
$data[44,3] = 44167
$data[44,4] = 9
$data[44,5] = 100112037
$data[44,6] = 'Cebollín'
$data[44,7] = 'Sin especificar'
$data[44,8] = 'Primera'
$data[44,9] = 65
$data[44,10] = 8000
$data[44,11] = 8000
$data[44,12] = 8000
$data[44,13] = '$/docena de paquetes'
$data[44,14] = 'Provincia de Cautín'
$data[44,15] = 667
$data[44,16] = 12
$data[44,17] = 'Hortaliza'

$data[45,0] = 10
$data[45,1] = 'Vega Modelo de Temuco'
$data[45,2] = 'La Araucanía'
$data[45,3] = 44167
$data[45,4] = 9
$data[45,5] = 100112037
$data[45,6] = 'Cebollín'
$data[45,7] = 'Sin especificar'
$data[45,8] = 'Primera'
$data[45,9] = 80
$data[45,10] = 7000
$data[45,11] = 7000
$data[45,12] = 7000
$data[45,13] = '$/media docena de atados'
$data[45,14] = 'Provincia de Cautín'
$data[45,15] = 1167
$data[45,16] = 6
$data[45,17] = 'Hortaliza'

$data[46,0] = 10
$data[46,1] = 'Vega Modelo de Temuco'
$data[46,2] = 'La Araucanía'
$data[46,3] = 44344
$data[46,4] = 9
$data[46,5] = 100112037
$data[46,6] = 'Cebollín'
$data[46,7] = 'Sin especificar'
$data[46,8] = 'Primera'
$data[46,9] = 50
$data[46,10] = 9000
$data[46,11] = 10000
$data[46,12] = 9400
$data[46,13] = '$/docena de paquetes'
$data[46,14] = 'Provincia de Cautín'
$data[46,15] = 783
$data[46,16] = 12
$data[46,17] = 'Hortaliza'

$data[47,0] = 10
$data[47,1] = 'Vega Modelo de Temuco'
$data[47,2] = 'La Araucanía'
$data[47,3] = 44217
$data[47,4] = 9
$data[47,5] = 100112037
$data[47,6] = 'Cebollín'
$data[47,7] = 'Sin especificar'
$data[47,8] = 'Primera'
$data[47,9] = 100
$data[47,10] = 5000
$data[47,11] = 5000
$data[47,12] = 5000
$data[47,13] = '$/docena de paquetes'
$data[47,14] = 'Provincia de Cautín'
$data[47,15] = 417
$data[47,16] = 12
$data[47,17] = 'Hortaliza'

$data[48,0] = 10
$data[48,1] = 'Vega Modelo de Temuco'
$data[48,2] = 'La Araucanía'
$data[48,3] = 44266
$data[48,4] = 9
$data[48,5] = 100112037
$data[48,6] = 'Cebollín'
$data[48,7] = 'Sin especificar'
$data[48,8] = 'Primera'
$data[48,9] = 75
$data[48,10] = 6000
$data[48,11] = 6000
$data[48,12] = 6000
$data[48,13] = '$/docena de paquetes'
$data[48,14] = 'Provincia de Cautín'
$data[48,15] = 500
$data[48,16] = 12
$data[48,17] = 'Hortaliza'

$data[49,0] = 10
$data[49,1] = 'Vega Modelo de Temuco'
$data[49,2] = 'La Araucanía'
$data[49,3] = 44350
$data[49,4] = 9
$data[49,5] = 100112037
$data[49,6] = 'Cebollín'
$data[49,7] = 'Sin especificar'
$data[49,8] = 'Primera'
$data[49,9] = 65
$data[49,10] = 9000
$data[49,11] = 10000
$data[49,12] = 9462
$data[49,13] = '$/docena de paquetes'
$data[49,14] = 'Provincia de Cautín'
$data[49,15] = 788
$data[49,16] = 12
$data[49,17] = 'Hortaliza'

$data[50,0] = 10
$data[50,1] = 'Vega Modelo de Temuco'
$data[50,2] = 'La Araucanía'
$data[50,3] = 44350
$data[50,4] = 9
$data[50,5] = 100112037
$data[50,6] = 'Cebollín'
$data[50,7] = 'Sin especificar'
$data[50,8] = 'Primera'
$data[50,9] = 70
$data[50,10] = 5000
$data[50,11] = 5000
$data[50,12] = 5000
$data[50,13] = '$/docena de paquetes'
$data[50,14] = 'Región Metropolitana'
$data[50,15] = 417
$data[50,16] = 12
$data[50,17] = 'Hortaliza'

$data[51,0] = 10
$data[51,1] = 'Vega Modelo de Temuco'
$data[51,2] = 'La Araucanía'
$data[51,3] = 44455
$data[51,4] = 9
$data[51,5] = 100112037
$data[51,6] = 'Cebollín'
$data[51,7] = 'Sin especificar'
$data[51,8] = 'Primera'
$data[51,9] = 40
$data[51,10] = 7000
$data[51,11] = 8000
$data[51,12] = 7500
$data[51,13] = '$/docena de paquetes'
$data[51,14] = 'Provincia de Cautín'
$data[51,15] = 625
$data[51,16] = 12
$data[51,17] = 'Hortaliza'

$data[52,0] = 10
$data[52,1] = 'Vega Modelo de Temuco'
$data[52,2] = 'La Araucanía'
$data[52,3] = 44455
$data[52,4] = 9
$data[52,5] = 100112037
$data[52,6] = 'Cebollín'
$data[52,7] = 'Sin especificar'
$data[52,8] = 'Primera'
$data[52,9] = 10
$data[52,10] = 5000
$data[52,11] = 5000
$data[52,12] = 5000
$data[52,13] = '$/docena de paquetes'
$data[52,14] = 'Región de O''Higgins'
$data[52,15] = 417
$data[52,16] = 12
$data[52,17] = 'Hortaliza'

$data[53,0] = 10
$data[53,1] = 'Vega Modelo de Temuco'
$data[53,2] = 'La Araucanía'
$data[53,3] = 44418
$data[53,4] = 9
$data[53,5] = 100112037
$data[53,6] = 'Cebollín'
$data[53,7] = 'Sin especificar'
$data[53,8] = 'Primera'
$data[53,9] = 35
$data[53,10] = 5000
$data[53,11] = 5000
$data[53,12] = 5000
$data[53,13] = '$/docena de atados'
$data[53,14] = 'Región de O''Higgins'
$data[53,15] = 1667
$data[53,16] = 3
$data[53,17] = 'Hortaliza'

$data[54,0] = 10
$data[54,1] = 'Vega Modelo de Temuco'
$data[54,2] = 'La Araucanía'
$data[54,3] = 44460
$data[54,4] = 9
$data[54,5] = 100112037
$data[54,6] = 'Cebollín'
$data[54,7] = 'Sin especificar'
$data[54,8] = 'Primera'
$data[54,9] = 30
$data[54,10] = 8000
$data[54,11] = 8000
$data[54,12] = 8000
$data[54,13] = '$/docena de paquetes'
$data[54,14] = 'Provincia de Cautín'
$data[54,15] = 667
$data[54,16] = 12
$data[54,17] = 'Hortaliza'

$data[55,0] = 10
$data[55,1] = 'Vega Modelo de Temuco'
$data[55,2] = 'La Araucanía'
$data[55,3] = 44460
$data[55,4] = 9
$data[55,5] = 100112037
$data[55,6] = 'Cebollín'
$data[55,7] = 'Sin especificar'
$data[55,8] = 'Primera'
$data[55,9] = 30
$data[55,10] = 5000
$data[55,11] = 5000
$data[55,12] = 5000
$data[55,13] = '$/docena de paquetes'
$data[55,14] = 'Región de O''Higgins'
$data[55,15] = 417
$data[55,16] = 12
$data[55,17] = 'Hortaliza'

$data[56,0] = 10
$data[56,1] = 'Vega Modelo de Temuco'
$data[56,2] = 'La Araucanía'
$data[56,3] = 44414
$data[56,4] = 9
$data[56,5] = 100112037
$data[56,6] = 'Cebollín'
$data[56,7] = 'Sin especificar'
$data[56,8] = 'Primera'
$data[56,9] = 20
$data[56,10] = 8000
$data[56,11] = 8000
$data[56,12] = 8000
$data[56,13] = '$/docena de paquetes'
$data[56,14] = 'Provincia de Cautín'
$data[56,15] = 667
$data[56,16] = 12
$data[56,17] = 'Hortaliza'

$data[57,0] = 10
$data[57,1] = 'Vega Modelo de Temuco'
$data[57,2] = 'La Araucanía'
$data[57,3] = 44165
$data[57,4] = 9
$data[57,5] = 100112037
$data[57,6] = 'Cebollín'
$data[57,7] = 'Sin especificar'
$data[57,8] = 'Primera'
$data[57,9] = 125
$data[57,10] = 8000
$data[57,11] = 8000
$data[57,12] = 8000
$data[57,13] = '$/docena de paquetes'
$data[57,14] = 'Provincia de Cautín'
$data[57,15] = 667
$data[57,16] = 12
$data[57,17] = 'Hortaliza'

$data[58,0] = 10
$data[58,1] = 'Vega Modelo de Temuco'
$data[58,2] = 'La Araucanía'
$data[58,3] = 44427
$data[58,4] = 9
$data[58,5] = 100112037
$data[58,6] = 'Cebollín'
$data[58,7] = 'Sin especificar'
$data[58,8] = 'Primera'
$data[58,9] = 40
$data[58,10] = 7000
$data[58,11] = 7000
$data[58,12] = 7000
$data[58,13] = '$/docena de paquetes'
$data[58,14] = 'Provincia de Cautín'
$data[58,15] = 583
$data[58,16] = 12
$data[58,17] = 'Hortaliza'

$data[59,0] = 10
$data[59,1] = 'Vega Modelo de Temuco'
$data[59,2] = 'La Araucanía'
$data[59,3] = 44427
$data[59,4] = 9
$data[59,5] = 100112037
$data[59,6] = 'Cebollín'
$data[59,7] = 'Sin especificar'
$data[59,8] = 'Primera'
$data[59,9] = 50
$data[59,10] = 5000
$data[59,11] = 5000
$data[59,12] = 5000
$data[59,13] = '$/docena de paquetes'
$data[59,14] = 'Región de O''Higgins'
$data[59,15] = 417
$data[59,16] = 12
$data[59,17] = 'Hortaliza'

$data[60,0] = 10
$data[60,1] = 'Vega Modelo de Temuco'
$data[60,2] = 'La Araucanía'
$data[60,3] = 44172
$data[60,4] = 9
$data[60,5] = 100112037
$data[60,6] = 'Cebollín'
$data[60,7] = 'Sin especificar'
$data[60,8] = 'Primera'
$data[60,9] = 110
$data[60,10] = 7000
$data[60,11] = 7000
$data[60,12] = 7000
$data[60,13] = '$/docena de paquetes'
$data[60,14] = 'Provincia de Cautín'
$data[60,15] = 583
$data[60,16] = 12
$data[60,17] = 'Hortaliza'

$data[61,0] = 10
$data[61,1] = 'Vega Modelo de Temuco'
$data[61,2] = 'La Araucanía'
$data[61,3] = 44172
$data[61,4] = 9
$data[61,5] = 100112037
$data[61,6] = 'Cebollín'
$data[61,7] = 'Sin especificar'
$data[61,8] = 'Primera'
$data[61,9] = 95
$data[61,10] = 6000
$data[61,11] = 6000
$data[61,12] = 6000
$data[61,13] = '$/docena de paquetes'
$data[61,14] = 'Región de O''Higgins'
$data[61,15] = 500
$data[61,16] = 12
$data[61,17] = 'Hortaliza'

$data[62,0] = 10
$data[62,1] = 'Vega Modelo de Temuco'
$data[62,2] = 'La Araucanía'
$data[62,3] = 44466
$data[62,4] = 9
$data[62,5] = 100112037
$data[62,6] = 'Cebollín'
$data[62,7] = 'Sin especificar'
$data[62,8] = 'Primera'
$data[62,9] = 50
$data[62,10] = 5000
$data[62,11] = 5000
$data[62,12] = 5000
$data[62,13] = '$/docena de paquetes'
$data[62,14] = 'Región de O''Higgins'
$data[62,15] = 417
$data[62,16] = 12
$data[62,17] = 'Hortaliza'

$data[63,0] = 10
$data[63,1] = 'Vega Modelo de Temuco'
$data[63,2] = 'La Araucanía'
$data[63,3] = 44389
$data[63,4] = 9
$data[63,5] = 100112037
$data[63,6] = 'Cebollín'
$data[63,7] = 'Sin especificar'
$data[63,8] = 'Primera'
$data[63,9] = 50
$data[63,10] = 8000
$data[63,11] = 8000
$data[63,12] = 8000
$data[63,13] = '$/docena de paquetes'
$data[63,14] = 'Provincia de Cautín'
$data[63,15] = 667
$data[63,16] = 12
$data[63,17] = 'Hortaliza'

$data[64,0] = 10
$data[64,1] = 'Vega Modelo de Temuco'
$data[64,2] = 'La Araucanía'
$data[64,3] = 44389
$data[64,4] = 9
$data[64,5] = 100112037
$data[64,6] = 'Cebollín'
$data[64,7] = 'Sin especificar'
$data[64,8] = 'Primera'
$data[64,9] = 155
$data[64,10] = 5000
$data[64,11] = 5000
$data[64,12] = 5000
$data[64,13] = '$/docena de paquetes'
$data[64,14] = 'Región de O''Higgins'
$data[64,15] = 417
$data[64,16] = 12
$data[64,17] = 'Hortaliza'

$data[65,0] = 10
$data[65,1] = 'Vega Modelo de Temuco'
$data[65,2] = 'La Araucanía'
$data[65,3] = 44249
$data[65,4] = 9
$data[65,5] = 100112037
$data[65,6] = 'Cebollín'
$data[65,7] = 'Sin especificar'
$data[65,8] = 'Primera'
$data[65,9] = 45
$data[65,10] = 6000
$data[65,11] = 6000
$data[65,12] = 6000
$data[65,13] = '$/docena de paquetes'
$data[65,14] = 'Provincia de Cautín'
$data[65,15] = 500
$data[65,16] = 12
$data[65,17] = 'Hortaliza'

$data[66,0] = 10
$data[66,1] = 'Vega Modelo de Temuco'
$data[66,2] = 'La Araucanía'
$data[66,3] = 44343
$data[66,4] = 9
$data[66,5] = 100112037
$data[66,6] = 'Cebollín'
$data[66,7] = 'Sin especificar'
$data[66,8] = 'Primera'
$data[66,9] = 60
$data[66,10] = 8000
$data[66,11] = 8000
$data[66,12] = 8000
$data[66,13] = '$/docena de paquetes'
$data[66,14] = 'Provincia de Cautín'
$data[66,15] = 667
$data[66,16] = 12
$data[66,17] = 'Hortaliza'

$data[67,0] = 10
$data[67,1] = 'Vega Modelo de Temuco'
$data[67,2] = 'La Araucanía'
$data[67,3] = 44201
$data[67,4] = 9
$data[67,5] = 100112037
$data[67,6] = 'Cebollín'
$data[67,7] = 'Sin especificar'
$data[67,8] = 'Primera'
$data[67,9] = 50
$data[67,10] = 7000
$data[67,11] = 8000
$data[67,12] = 7400
$data[67,13] = '$/docena de paquetes'
$data[67,14] = 'Provincia de Cautín'
$data[67,15] = 617
$data[67,16] = 12
$data[67,17] = 'Hortaliza'

$data[68,0] = 10
$data[68,1] = 'Vega Modelo de Temuco'
$data[68,2] = 'La Araucanía'
$data[68,3] = 44280
$data[68,4] = 9
$data[68,5] = 100112037
$data[68,6] = 'Cebollín'
$data[68,7] = 'Sin especificar'
$data[68,8] = 'Primera'
$data[68,9] = 95
$data[68,10] = 5000
$data[68,11] = 5000
$data[68,12] = 5000
$data[68,13] = '$/docena de paquetes'
$data[68,14] = 'Provincia de Cautín'
$data[68,15] = 417
$data[68,16] = 12
$data[68,17] = 'Hortaliza'

$data[69,0] = 10
$data[69,1] = 'Vega Modelo de Temuco'
$data[69,2] = 'La Araucanía'
$data[69,3] = 44447
$data[69,4] = 9
$data[69,5] = 100112037
$data[69,6] = 'Cebollín'
$data[69,7] = 'Sin especificar'
$data[69,8] = 'Primera'
$data[69,9] = 55
$data[69,10] = 5000
$data[69,11] = 5000
$data[69,12] = 5000
$data[69,13] = '$/docena de paquetes'
$data[69,14] = 'Región de O''Higgins'
$data[69,15] = 417
$data[69,16] = 12
$data[69,17] = 'Hortaliza'

$data[70,0] = 10
$data[70,1] = 'Vega Modelo de Temuco'
$data[70,2] = 'La Araucanía'
$data[70,3] = 44270
$data[70,4] = 9
$data[70,5] = 100112037
$data[70,6] = 'Cebollín'
$data[70,7] = 'Sin especificar'
$data[70,8] = 'Primera'
$data[70,9] = 70
$data[70,10] = 5000
$data[70,11] = 5000
$data[70,12] = 5000
$data[70,13] = '$/docena de paquetes'
$data[70,14] = 'Provincia de Cautín'
$data[70,15] = 417
$data[70,16] = 12
$data[70,17] = 'Hortaliza'

$data[71,0] = 10
$data[71,1] = 'Vega Modelo de Temuco'
$data[71,2] = 'La Araucanía'
$data[71,3] = 44312
$data[71,4] = 9
$data[71,5] = 100112037
$data[71,6] = 'Cebollín'
$data[71,7] = 'Sin especificar'
$data[71,8] = 'Primera'
$data[71,9] = 70
$data[71,10] = 7000
$data[71,11] = 7000
$data[71,12] = 7000
$data[71,13] = '$/docena de paquetes'
$data[71,14] = 'Provincia de Cautín'
$data[71,15] = 583
$data[71,16] = 12
$data[71,17] = 'Hortaliza'

$data[72,0] = 10
$data[72,1] = 'Vega Modelo de Temuco'
$data[72,2] = 'La Araucanía'
$data[72,3] = 44187
$data[72,4] = 9
$data[72,5] = 100112037
$data[72,6] = 'Cebollín'
$data[72,7] = 'Sin especificar'
$data[72,8] = 'Primera'
$data[72,9] = 70
$data[72,10] = 7000
$data[72,11] = 8000
$data[72,12] = 7429
$data[72,13] = '$/docena de paquetes'
$data[72,14] = 'Provincia de Cautín'
$data[72,15] = 619
$data[72,16] = 12
$data[72,17] = 'Hortaliza'

$data[73,0] = 10
$data[73,1] = 'Vega Modelo de Temuco'
$data[73,2] = 'La Araucanía'
$data[73,3] = 44390
$data[73,4] = 9
$data[73,5] = 100112037
$data[73,6] = 'Cebollín'
$data[73,7] = 'Sin especificar'
$data[73,8] = 'Primera'
$data[73,9] = 120
$data[73,10] = 9000
$data[73,11] = 10000
$data[73,12] = 9458
$data[73,13] = '$/docena de paquetes'
$data[73,14] = 'Provincia de Cautín'
$data[73,15] = 788
$data[73,16] = 12
$data[73,17] = 'Hortaliza'

$data[74,0] = 10
$data[74,1] = 'Vega Modelo de Temuco'
$data[74,2] = 'La Araucanía'
$data[74,3] = 44386
$data[74,4] = 9
$data[74,5] = 100112037
$data[74,6] = 'Cebollín'
$data[74,7] = 'Sin especificar'
$data[74,8] = 'Primera'
$data[74,9] = 20
$data[74,10] = 8000
$data[74,11] = 8000
$data[74,12] = 8000
$data[74,13] = '$/docena de paquetes'
$data[74,14] = 'Provincia de Cautín'
$data[74,15] = 667
$data[74,16] = 12
$data[74,17] = 'Hortaliza'

$data[75,0] = 10
$data[75,1] = 'Vega Modelo de Temuco'
$data[75,2] = 'La Araucanía'
$data[75,3] = 44308
$data[75,4] = 9
$data[75,5] = 100112037
$data[75,6] = 'Cebollín'
$data[75,7] = 'Sin especificar'
$data[75,8] = 'Primera'
$data[75,9] = 80
$data[75,10] = 6000
$data[75,11] = 7000
$data[75,12] = 6438
$data[75,13] = '$/docena de paquetes'
$data[75,14] = 'Provincia de Cautín'
$data[75,15] = 536
$data[75,16] = 12
$data[75,17] = 'Hortaliza'

$data[76,0] = 10
$data[76,1] = 'Vega Modelo de Temuco'
$data[76,2] = 'La Araucanía'
$data[76,3] = 44463
$data[76,4] = 9
$data[76,5] = 100112037
$data[76,6] = 'Cebollín'
$data[76,7] = 'Sin especificar'
$data[76,8] = 'Primera'
$data[76,9] = 20
$data[76,10] = 8000
$data[76,11] = 8000
$data[76,12] = 8000
$data[76,13] = '$/docena de paquetes'
$data[76,14] = 'Provincia de Cautín'
$data[76,15] = 667
$data[76,16] = 12
$data[76,17] = 'Hortaliza'

$data[77,0] = 10
$data[77,1] = 'Vega Modelo de Temuco'
$data[77,2] = 'La Araucanía'
$data[77,3] = 44196
$data[77,4] = 9
$data[77,5] = 100112037
$data[77,6] = 'Cebollín'
$data[77,7] = 'Sin especificar'
$data[77,8] = 'Primera'
$data[77,9] = 30
$data[77,10] = 5000
$data[77,11] = 5000
$data[77,12] = 5000
$data[77,13] = '$/docena de paquetes'
$data[77,14] = 'Provincia de Cautín'
$data[77,15] = 417
$data[77,16] = 12
$data[77,17] = 'Hortaliza'

$data[78,0] = 10
$data[78,1] = 'Vega Modelo de Temuco'
$data[78,2] = 'La Araucanía'
$data[78,3] = 44251
$data[78,4] = 9
$data[78,5] = 100112037
$data[78,6] = 'Cebollín'
$data[78,7] = 'Sin especificar'
$data[78,8] = 'Primera'
$data[78,9] = 85
$data[78,10] = 5000
$data[78,11] = 6000
$data[78,12] = 5529
$data[78,13] = '$/docena de paquetes'
$data[78,14] = 'Provincia de Cautín'
$data[78,15] = 461
$data[78,16] = 12
$data[78,17] = 'Hortaliza'

$data[79,0] = 10
$data[79,1] = 'Vega Modelo de Temuco'
$data[79,2] = 'La Araucanía'
$data[79,3] = 44243
$data[79,4] = 9
$data[79,5] = 100112037
$data[79,6] = 'Cebollín'
$data[79,7] = 'Sin especificar'
$data[79,8] = 'Primera'
$data[79,9] = 65
$data[79,10] = 7000
$data[79,11] = 7000
$data[79,12] = 7000
$data[79,13] = '$/docena de paquetes'
$data[79,14] = 'Provincia de Cautín'
$data[79,15] = 583
$data[79,16] = 12
$data[79,17] = 'Hortaliza'

$data[80,0] = 10
$data[80,1] = 'Vega Modelo de Temuco'
$data[80,2] = 'La Araucanía'
$data[80,3] = 44252
$data[80,4] = 9
$data[80,5] = 100112037
$data[80,6] = 'Cebollín'
$data[80,7] = 'Sin especificar'
$data[80,8] = 'Primera'
$data[80,9] = 140
$data[80,10] = 5000
$data[80,11] = 6000
$data[80,12] = 5393
$data[80,13] = '$/docena de paquetes'
$data[80,14] = 'Provincia de Cautín'
$data[80,15] = 449
$data[80,16] = 12
$data[80,17] = 'Hortaliza'

$data[81,0] = 10
$data[81,1] = 'Vega Modelo de Temuco'
$data[81,2] = 'La Araucanía'
$data[81,3] = 44166
$data[81,4] = 9
$data[81,5] = 100112037
$data[81,6] = 'Cebollín'
$data[81,7] = 'Sin especificar'
$data[81,8] = 'Primera'
$data[81,9] = 75
$data[81,10] = 8000
$data[81,11] = 8000
$data[81,12] = 8000
$data[81,13] = '$/docena de paquetes'
$data[81,14] = 'Provincia de Cautín'
$data[81,15] = 667
$data[81,16] = 12
$data[81,17] = 'Hortaliza'

$data[82,0] = 10
$data[82,1] = 'Vega Modelo de Temuco'
$data[82,2] = 'La Araucanía'
$data[82,3] = 44168
$data[82,4] = 9
$data[82,5] = 100112037
$data[82,6] = 'Cebollín'
$data[82,7] = 'Sin especificar'
$data[82,8] = 'Primera'
$data[82,9] = 150
$data[82,10] = 7000
$data[82,11] = 8000
$data[82,12] = 7433
$data[82,13] = '$/docena de paquetes'
$data[82,14] = 'Provincia de Cautín'
$data[82,15] = 619
$data[82,16] = 12
$data[82,17] = 'Hortaliza'

$data[83,0] = 10
$data[83,1] = 'Vega Modelo de Temuco'
$data[83,2] = 'La Araucanía'
$data[83,3] = 44369
$data[83,4] = 9
$data[83,5] = 100112037
$data[83,6] = 'Cebollín'
$data[83,7] = 'Sin especificar'
$data[83,8] = 'Primera'
$data[83,9] = 30
$data[83,10] = 8000
$data[83,11] = 8000
$data[83,12] = 8000
$data[83,13] = '$/docena de paquetes'
$data[83,14] = 'Provincia de Cautín'
$data[83,15] = 667
$data[83,16] = 12
$data[83,17] = 'Hortaliza'

$data[84,0] = 10
$data[84,1] = 'Vega Modelo de Temuco'
$data[84,2] = 'La Araucanía'
$data[84,3] = 44369
$data[84,4] = 9
$data[84,5] = 100112037
$data[84,6] = 'Cebollín'
$data[84,7] = 'Sin especificar'
$data[84,8] = 'Primera'
$data[84,9] = 30
$data[84,10] = 5000
$data[84,11] = 5000
$data[84,12] = 5000
$data[84,13] = '$/docena de paquetes'
$data[84,14] = 'Región Metropolitana'
$data[84,15] = 417
$data[84,16] = 12
$data[84,17] = 'Hortaliza'

$data[85,0] = 10
$data[85,1] = 'Vega Modelo de Temuco'
$data[85,2] = 'La Araucanía'
$data[85,3] = 44433
$data[85,4] = 9
$data[85,5] = 100112037
$data[85,6] = 'Cebollín'
$data[85,7] = 'Sin especificar'
$data[85,8] = 'Primera'
$data[85,9] = 35
$data[85,10] = 5000
$data[85,11] = 5000
$data[85,12] = 5000
$data[85,13] = '$/docena de paquetes'
$data[85,14] = 'Región de O''Higgins'
$data[85,15] = 417
$data[85,16] = 12
$data[85,17] = 'Hortaliza'

$data[86,0] = 10
$data[86,1] = 'Vega Modelo de Temuco'
$data[86,2] = 'La Araucanía'
$data[86,3] = 44221
$data[86,4] = 9
$data[86,5] = 100112037
$data[86,6] = 'Cebollín'
$data[86,7] = 'Sin especificar'
$data[86,8] = 'Primera'
$data[86,9] = 95
$data[86,10] = 7000
$data[86,11] = 7000
$data[86,12] = 7000
$data[86,13] = '$/docena de paquetes'
$data[86,14] = 'Provincia de Cautín'
$data[86,15] = 583
$data[86,16] = 12
$data[86,17] = 'Hortaliza'

$data[87,0] = 10
$data[87,1] = 'Vega Modelo de Temuco'
$data[87,2] = 'La Araucanía'
$data[87,3] = 44371
$data[87,4] = 9
$data[87,5] = 100112037
$data[87,6] = 'Cebollín'
$data[87,7] = 'Sin especificar'
$data[87,8] = 'Primera'
$data[87,9] = 20
$data[87,10] = 8000
$data[87,11] = 8000
$data[87,12] = 8000
$data[87,13] = '$/docena de paquetes'
$data[87,14] = 'Provincia de Cautín'
$data[87,15] = 667
$data[87,16] = 12
$data[87,17] = 'Hortaliza'

$data[88,0] = 10
$data[88,1] = 'Vega Modelo de Temuco'
$data[88,2] = 'La Araucanía'
$data[88,3] = 44371
$data[88,4] = 9
$data[88,5] = 100112037
$data[88,6] = 'Cebollín'
$data[88,7] = 'Sin especificar'
$data[88,8] = 'Primera'
$data[88,9] = 30
$data[88,10] = 5000
$data[88,11] = 5000
$data[88,12] = 5000
$data[88,13] = '$/docena de paquetes'
$data[88,14] = 'Región Metropolitana'
$data[88,15] = 417
$data[88,16] = 12
$data[88,17] = 'Hortaliza'

$data[89,0] = 10
$data[89,1] = 'Vega Modelo de Temuco'
$data[89,2] = 'La Araucanía'
$data[89,3] = 44371
$data[89,4] = 9
$data[89,5] = 100112037
$data[89,6] = 'Cebollín'
$data[89,7] = 'Sin especificar'
$data[89,8] = 'Primera'
$data[89,9] = 40
$data[89,10] = 5000
$data[89,11] = 5000
$data[89,12] = 5000
$data[89,13] = '$/docena de paquetes'
$data[89,14] = 'Región del Maule'
$data[89,15] = 417
$data[89,16] = 12
$data[89,17] = 'Hortaliza'

$data[90,0] = 10
$data[90,1] = 'Vega Modelo de Temuco'
$data[90,2] = 'La Araucanía'
$data[90,3] = 44316
$data[90,4] = 9
$data[90,5] = 100112037
$data[90,6] = 'Cebollín'
$data[90,7] = 'Sin especificar'
$data[90,8] = 'Primera'
$data[90,9] = 40
$data[90,10] = 7000
$data[90,11] = 7000
$data[90,12] = 7000
$data[90,13] = '$/docena de paquetes'
$data[90,14] = 'Provincia de Cautín'
$data[90,15] = 583
$data[90,16] = 12
$data[90,17] = 'Hortaliza'

$data[91,0] = 10
$data[91,1] = 'Vega Modelo de Temuco'
$data[91,2] = 'La Araucanía'
$data[91,3] = 44397
$data[91,4] = 9
$data[91,5] = 100112037
$data[91,6] = 'Cebollín'
$data[91,7] = 'Sin especificar'
$data[91,8] = 'Primera'
$data[91,9] = 70
$data[91,10] = 9000
$data[91,11] = 10000
$data[91,12] = 9429
$data[91,13] = '$/docena de paquetes'
$data[91,14] = 'Provincia de Cautín'
$data[91,15] = 786
$data[91,16] = 12
$data[91,17] = 'Hortaliza'

$data[92,0] = 10
$data[92,1] = 'Vega Modelo de Temuco'
$data[92,2] = 'La Araucanía'
$data[92,3] = 44397
$data[92,4] = 9
$data[92,5] = 100112037
$data[92,6] = 'Cebollín'
$data[92,7] = 'Sin especificar'
$data[92,8] = 'Primera'
$data[92,9] = 30
$data[92,10] = 5000
$data[92,11] = 5000
$data[92,12] = 5000
$data[92,13] = '$/docena de paquetes'
$data[92,14] = 'Región de O''Higgins'
$data[92,15] = 417
$data[92,16] = 12
$data[92,17] = 'Hortaliza'

$data[93,0] = 10
$data[93,1] = 'Vega Modelo de Temuco'
$data[93,2] = 'La Araucanía'
$data[93,3] = 44363
$data[93,4] = 9
$data[93,5] = 100112037
$data[93,6] = 'Cebollín'
$data[93,7] = 'Sin especificar'
$data[93,8] = 'Primera'
$data[93,9] = 75
$data[93,10] = 5000
$data[93,11] = 5000
$data[93,12] = 5000
$data[93,13] = '$/docena de atados'
$data[93,14] = 'Región Metropolitana'
$data[93,15] = 1667
$data[93,16] = 3
$data[93,17] = 'Hortaliza'

$data[94,0] = 10
$data[94,1] = 'Vega Modelo de Temuco'
$data[94,2] = 'La Araucanía'
$data[94,3] = 44277
$data[94,4] = 9
$data[94,5] = 100112037
$data[94,6] = 'Cebollín'
$data[94,7] = 'Sin especificar'
$data[94,8] = 'Primera'
$data[94,9] = 75
$data[94,10] = 5000
$data[94,11] = 6000
$data[94,12] = 5467
$data[94,13] = '$/docena de atados'
$data[94,14] = 'Provincia de Cautín'
$data[94,15] = 1822
$data[94,16] = 3
$data[94,17] = 'Hortaliza'

$data[95,0] = 10
$data[95,1] = 'Vega Modelo de Temuco'
$data[95,2] = 'La Araucanía'
$data[95,3] = 44291
$data[95,4] = 9
$data[95,5] = 100112037
$data[95,6] = 'Cebollín'
$data[95,7] = 'Sin especificar'
$data[95,8] = 'Primera'
$data[95,9] = 55
$data[95,10] = 7000
$data[95,11] = 7000
$data[95,12] = 7000
$data[95,13] = '$/docena de paquetes'
$data[95,14] = 'Provincia de Cautín'
$data[95,15] = 583
$data[95,16] = 12
$data[95,17] = 'Hortaliza'

$data[96,0] = 10
$data[96,1] = 'Vega Modelo de Temuco'
$data[96,2] = 'La Araucanía'
$data[96,3] = 44273
$data[96,4] = 9
$data[96,5] = 100112037
$data[96,6] = 'Cebollín'
$data[96,7] = 'Sin especificar'
$data[96,8] = 'Primera'
$data[96,9] = 50
$data[96,10] = 5000
$data[96,11] = 5000
$data[96,12] = 5000
$data[96,13] = '$/docena de paquetes'
$data[96,14] = 'Provincia de Cautín'
$data[96,15] = 417
$data[96,16] = 12
$data[96,17] = 'Hortaliza'

$data[97,0] = 10
$data[97,1] = 'Vega Modelo de Temuco'
$data[97,2] = 'La Araucanía'
$data[97,3] = 44438
$data[97,4] = 9
$data[97,5] = 100112037
$data[97,6] = 'Cebollín'
$data[97,7] = 'Sin especificar'
$data[97,8] = 'Primera'
$data[97,9] = 30
$data[97,10] = 9000
$data[97,11] = 9000
$data[97,12] = 9000
$data[97,13] = '$/docena de paquetes'
$data[97,14] = 'Provincia de Cautín'
$data[97,15] = 750
$data[97,16] = 12
$data[97,17] = 'Hortaliza'

$data[98,0] = 10
$data[98,1] = 'Vega Modelo de Temuco'
$data[98,2] = 'La Araucanía'
$data[98,3] = 44372
$data[98,4] = 9
$data[98,5] = 100112037
$data[98,6] = 'Cebollín'
$data[98,7] = 'Sin especificar'
$data[98,8] = 'Primera'
$data[98,9] = 30
$data[98,10] = 8000
$data[98,11] = 8000
$data[98,12] = 8000
$data[98,13] = '$/docena de paquetes'
$data[98,14] = 'Provincia de Cautín'
$data[98,15] = 667
$data[98,16] = 12
$data[98,17] = 'Hortaliza'

$data[99,0] = 10
$data[99,1] = 'Vega Modelo de Temuco'
$data[99,2] = 'La Araucanía'
$data[99,3] = 44286
$data[99,4] = 9
$data[99,5] = 100112037
$data[99,6] = 'Cebollín'
$data[99,7] = 'Sin especificar'
$data[99,8] = 'Primera'
$data[99,9] = 30
$data[99,10] = 7000
$data[99,11] = 7000
$data[99,12] = 7000
$data[99,13] = '$/docena de paquetes'
$data[99,14] = 'Provincia de Cautín'
$data[99,15] = 583
$data[99,16] = 12
$data[99,17] = 'Hortaliza'

$data[100,0] = 10
$data[100,1] = 'Vega Modelo de Temuco'
$data[100,2] = 'La Araucanía'
$data[100,3] = 44209
$data[100,4] = 9
$data[100,5] = 100112037
$data[100,6] = 'Cebollín'
$data[100,7] = 'Sin especificar'
$data[100,8] = 'Primera'
$data[100,9] = 65
$data[100,10] = 5000
$data[100,11] = 5000
$data[100,12] = 5000
$data[100,13] = '$/docena de paquetes'
$data[100,14] = 'Región de O''Higgins'
$data[100,15] = 417
$data[100,16] = 12
$data[100,17] = 'Hortaliza'

$data[101,0] = 10
$data[101,1] = 'Vega Modelo de Temuco'
$data[101,2] = 'La Araucanía'
$data[101,3] = 44356
$data[101,4] = 9
$data[101,5] = 100112037
$data[101,6] = 'Cebollín'
$data[101,7] = 'Sin especificar'
$data[101,8] = 'Primera'
$data[101,9] = 20
$data[101,10] = 9000
$data[101,11] = 9000
$data[101,12] = 9000
$data[101,13] = '$/docena de paquetes'
$data[101,14] = 'Provincia de Cautín'
$data[101,15] = 750
$data[101,16] = 12
$data[101,17] = 'Hortaliza'

$data[102,0] = 10
$data[102,1] = 'Vega Modelo de Temuco'
$data[102,2] = 'La Araucanía'
$data[102,3] = 44356
$data[102,4] = 9
$data[102,5] = 100112037
$data[102,6] = 'Cebollín'
$data[102,7] = 'Sin especificar'
$data[102,8] = 'Primera'
$data[102,9] = 30
$data[102,10] = 5000
$data[102,11] = 5000
$data[102,12] = 5000
$data[102,13] = '$/docena de paquetes'
$data[102,14] = 'Región Metropolitana'
$data[102,15] = 417
$data[102,16] = 12
$data[102,17] = 'Hortaliza'

$data[103,0] = 10
$data[103,1] = 'Vega Modelo de Temuco'
$data[103,2] = 'La Araucanía'
$data[103,3] = 44160
$data[103,4] = 9
$data[103,5] = 100112037
$data[103,6] = 'Cebollín'
$data[103,7] = 'Sin especificar'
$data[103,8] = 'Primera'
$data[103,9] = 50
$data[103,10] = 5000
$data[103,11] = 5000
$data[103,12] = 5000
$data[103,13] = '$/docena de paquetes'
$data[103,14] = 'Región de O''Higgins'
$data[103,15] = 417
$data[103,16] = 12
$data[103,17] = 'Hortaliza'

$data[104,0] = 10
$data[104,1] = 'Vega Modelo de Temuco'
$data[104,2] = 'La Araucanía'
$data[104,3] = 44351
$data[104,4] = 9
$data[104,5] = 100112037
$data[104,6] = 'Cebollín'
$data[104,7] = 'Sin especificar'
$data[104,8] = 'Primera'
$data[104,9] = 135
$data[104,10] = 9000
$data[104,11] = 9000
$data[104,12] = 9000
$data[104,13] = '$/docena de paquetes'
$data[104,14] = 'Provincia de Cautín'
$data[104,15] = 750
$data[104,16] = 12
$data[104,17] = 'Hortaliza'

$data[105,0] = 10
$data[105,1] = 'Vega Modelo de Temuco'
$data[105,2] = 'La Araucanía'
$data[105,3] = 44351
$data[105,4] = 9
$data[105,5] = 100112037
$data[105,6] = 'Cebollín'
$data[105,7] = 'Sin especificar'
$data[105,8] = 'Primera'
$data[105,9] = 65
$data[105,10] = 5000
$data[105,11] = 5000
$data[105,12] = 5000
$data[105,13] = '$/docena de paquetes'
$data[105,14] = 'Región Metropolitana'
$data[105,15] = 417
$data[105,16] = 12
$data[105,17] = 'Hortaliza'

$data[106,0] = 10
$data[106,1] = 'Vega Modelo de Temuco'
$data[106,2] = 'La Araucanía'
$data[106,3] = 44365
$data[106,4] = 9
$data[106,5] = 100112037
$data[106,6] = 'Cebollín'
$data[106,7] = 'Sin especificar'
$data[106,8] = 'Primera'
$data[106,9] = 95
$data[106,10] = 8000
$data[106,11] = 8000
$data[106,12] = 8000
$data[106,13] = '$/docena de paquetes'
$data[106,14] = 'Provincia de Cautín'
$data[106,15] = 667
$data[106,16] = 12
$data[106,17] = 'Hortaliza'

$data[107,0] = 10
$data[107,1] = 'Vega Modelo de Temuco'
$data[107,2] = 'La Araucanía'
$data[107,3] = 44215
$data[107,4] = 9
$data[107,5] = 100112037
$data[107,6] = 'Cebollín'
$data[107,7] = 'Sin especificar'
$data[107,8] = 'Primera'
$data[107,9] = 160
$data[107,10] = 5000
$data[107,11] = 6000
$data[107,12] = 5500
$data[107,13] = '$/docena de paquetes'
$data[107,14] = 'Provincia de Cautín'
$data[107,15] = 458
$data[107,16] = 12
$data[107,17] = 'Hortaliza'

$data[108,0] = 10
$data[108,1] = 'Vega Modelo de Temuco'
$data[108,2] = 'La Araucanía'
$data[108,3] = 44175
$data[108,4] = 9
$data[108,5] = 100112037
$data[108,6] = 'Cebollín'
$data[108,7] = 'Sin especificar'
$data[108,8] = 'Primera'
$data[108,9] = 50
$data[108,10] = 8000
$data[108,11] = 8000
$data[108,12] = 8000
$data[108,13] = '$/docena de paquetes'
$data[108,14] = 'Provincia de Cautín'
$data[108,15] = 667
$data[108,16] = 12
$data[108,17] = 'Hortaliza'

$data[109,0] = 10
$data[109,1] = 'Vega Modelo de Temuco'
$data[109,2] = 'La Araucanía'
$data[109,3] = 44461
$data[109,4] = 9
$data[109,5] = 100112037
$data[109,6] = 'Cebollín'
$data[109,7] = 'Sin especificar'
$data[109,8] = 'Primera'
$data[109,9] = 20
$data[109,10] = 8000
$data[109,11] = 8000
$data[109,12] = 8000
$data[109,13] = '$/docena de paquetes'
$data[109,14] = 'Provincia de Cautín'
$data[109,15] = 667
$data[109,16] = 12
$data[109,17] = 'Hortaliza'

$data[110,0] = 10
$data[110,1] = 'Vega Modelo de Temuco'
$data[110,2] = 'La Araucanía'
$data[110,3] = 44461
$data[110,4] = 9
$data[110,5] = 100112037
$data[110,6] = 'Cebollín'
$data[110,7] = 'Sin especificar'
$data[110,8] = 'Primera'
$data[110,9] = 40
$data[110,10] = 5000
$data[110,11] = 5000
$data[110,12] = 5000
$data[110,13] = '$/docena de paquetes'
$data[110,14] = 'Región de O''Higgins'
$data[110,15] = 417
$data[110,16] = 12
$data[110,17] = 'Hortaliza'

$data[111,0] = 10
$data[111,1] = 'Vega Modelo de Temuco'
$data[111,2] = 'La Araucanía'
$data[111,3] = 44357
$data[111,4] = 9
$data[111,5] = 100112037
$data[111,6] = 'Cebollín'
$data[111,7] = 'Sin especificar'
$data[111,8] = 'Primera'
$data[111,9] = 80
$data[111,10] = 9000
$data[111,11] = 9000
$data[111,12] = 9000
$data[111,13] = '$/docena de paquetes'
$data[111,14] = 'Provincia de Cautín'
$data[111,15] = 750
$data[111,16] = 12
$data[111,17] = 'Hortaliza'

$data[112,0] = 10
$data[112,1] = 'Vega Modelo de Temuco'
$data[112,2] = 'La Araucanía'
$data[112,3] = 44203
$data[112,4] = 9
$data[112,5] = 100112037
$data[112,6] = 'Cebollín'
$data[112,7] = 'Sin especificar'
$data[112,8] = 'Primera'
$data[112,9] = 100
$data[112,10] = 8000
$data[112,11] = 8000
$data[112,12] = 8000
$data[112,13] = '$/docena de paquetes'
$data[112,14] = 'Provincia de Cautín'
$data[112,15] = 667
$data[112,16] = 12
$data[112,17] = 'Hortaliza'

$data[113,0] = 10
$data[113,1] = 'Vega Modelo de Temuco'
$data[113,2] = 'La Araucanía'
$data[113,3] = 44162
$data[113,4] = 9
$data[113,5] = 100112037
$data[113,6] = 'Cebollín'
$data[113,7] = 'Sin especificar'
$data[113,8] = 'Primera'
$data[113,9] = 50
$data[113,10] = 7000
$data[113,11] = 7000
$data[113,12] = 7000
$data[113,13] = '$/docena de paquetes'
$data[113,14] = 'Provincia de Cautín'
$data[113,15] = 583
$data[113,16] = 12
$data[113,17] = 'Hortaliza'

$data[114,0] = 10
$data[114,1] = 'Vega Modelo de Temuco'
$data[114,2] = 'La Araucanía'
$data[114,3] = 44162
$data[114,4] = 9
$data[114,5] = 100112037
$data[114,6] = 'Cebollín'
$data[114,7] = 'Sin especificar'
$data[114,8] = 'Primera'
$data[114,9] = 40
$data[114,10] = 5000
$data[114,11] = 5000
$data[114,12] = 5000
$data[114,13] = '$/docena de paquetes'
$data[114,14] = 'Región del Maule'
$data[114,15] = 417
$data[114,16] = 12
$data[114,17] = 'Hortaliza'

$data[115,0] = 10
$data[115,1] = 'Vega Modelo de Temuco'
$data[115,2] = 'La Araucanía'
$data[115,3] = 44410
$data[115,4] = 9
$data[115,5] = 100112037
$data[115,6] = 'Cebollín'
$data[115,7] = 'Sin especificar'
$data[115,8] = 'Primera'
$data[115,9] = 60
$data[115,10] = 9000
$data[115,11] = 10000
$data[115,12] = 9500
$data[115,13] = '$/docena de paquetes'
$data[115,14] = 'Provincia de Cautín'
$data[115,15] = 792
$data[115,16] = 12
$data[115,17] = 'Hortaliza'

$data[116,0] = 10
$data[116,1] = 'Vega Modelo de Temuco'
$data[116,2] = 'La Araucanía'
$data[116,3] = 44410
$data[116,4] = 9
$data[116,5] = 100112037
$data[116,6] = 'Cebollín'
$data[116,7] = 'Sin especificar'
$data[116,8] = 'Primera'
$data[116,9] = 50
$data[116,10] = 5000
$data[116,11] = 5000
$data[116,12] = 5000
$data[116,13] = '$/docena de paquetes'
$data[116,14] = 'Región de O''Higgins'
$data[116,15] = 417
$data[116,16] = 12
$data[116,17] = 'Hortaliza'

$data[117,0] = 10
$data[117,1] = 'Vega Modelo de Temuco'
$data[117,2] = 'La Araucanía'
$data[117,3] = 44411
$data[117,4] = 9
$data[117,5] = 100112037
$data[117,6] = 'Cebollín'
$data[117,7] = 'Sin especificar'
$data[117,8] = 'Primera'
$data[117,9] = 30
$data[117,10] = 10000
$data[117,11] = 10000
$data[117,12] = 10000
$data[117,13] = '$/docena de paquetes'
$data[117,14] = 'Provincia de Cautín'
$data[117,15] = 833
$data[117,16] = 12
$data[117,17] = 'Hortaliza'

$data[118,0] = 10
$data[118,1] = 'Vega Modelo de Temuco'
$data[118,2] = 'La Araucanía'
$data[118,3] = 44176
$data[118,4] = 9
$data[118,5] = 100112037
$data[118,6] = 'Cebollín'
$data[118,7] = 'Sin especificar'
$data[118,8] = 'Primera'
$data[118,9] = 30
$data[118,10] = 8000
$data[118,11] = 9000
$data[118,12] = 8667
$data[118,13] = '$/docena de paquetes'
$data[118,14] = 'Provincia de Cautín'
$data[118,15] = 722
$data[118,16] = 12
$data[118,17] = 'Hortaliza'

$data[119,0] = 10
$data[119,1] = 'Vega Modelo de Temuco'
$data[119,2] = 'La Araucanía'
$data[119,3] = 44239
$data[119,4] = 9
$data[119,5] = 100112037
$data[119,6] = 'Cebollín'
$data[119,7] = 'Sin especificar'
$data[119,8] = 'Primera'
$data[119,9] = 115
$data[119,10] = 6000
$data[119,11] = 7000
$data[119,12] = 6565
$data[119,13] = '$/docena de paquetes'
$data[119,14] = 'Provincia de Cautín'
$data[119,15] = 547
$data[119,16] = 12
$data[119,17] = 'Hortaliza'

$data[120,0] = 10
$data[120,1] = 'Vega Modelo de Temuco'
$data[120,2] = 'La Araucanía'
$data[120,3] = 44376
$data[120,4] = 9
$data[120,5] = 100112037
$data[120,6] = 'Cebollín'
$data[120,7] = 'Sin especificar'
$data[120,8] = 'Primera'
$data[120,9] = 65
$data[120,10] = 5000
$data[120,11] = 5000
$data[120,12] = 5000
$data[120,13] = '$/docena de paquetes'
$data[120,14] = 'Región Metropolitana'
$data[120,15] = 417
$data[120,16] = 12
$data[120,17] = 'Hortaliza'

$data[121,0] = 10
$data[121,1] = 'Vega Modelo de Temuco'
$data[121,2] = 'La Araucanía'
$data[121,3] = 44292
$data[121,4] = 9
$data[121,5] = 100112037
$data[121,6] = 'Cebollín'
$data[121,7] = 'Sin especificar'
$data[121,8] = 'Primera'
$data[121,9] = 35
$data[121,10] = 7000
$data[121,11] = 7000
$data[121,12] = 7000
$data[121,13] = '$/docena de paquetes'
$data[121,14] = 'Provincia de Cautín'
$data[121,15] = 583
$data[121,16] = 12
$data[121,17] = 'Hortaliza'

$data[122,0] = 10
$data[122,1] = 'Vega Modelo de Temuco'
$data[122,2] = 'La Araucanía'
$data[122,3] = 44358
$data[122,4] = 9
$data[122,5] = 100112037
$data[122,6] = 'Cebollín'
$data[122,7] = 'Sin especificar'
$data[122,8] = 'Primera'
$data[122,9] = 30
$data[122,10] = 9000
$data[122,11] = 9000
$data[122,12] = 9000
$data[122,13] = '$/docena de paquetes'
$data[122,14] = 'Provincia de Cautín'
$data[122,15] = 750
$data[122,16] = 12
$data[122,17] = 'Hortaliza'

$data[123,0] = 10
$data[123,1] = 'Vega Modelo de Temuco'
$data[123,2] = 'La Araucanía'
$data[123,3] = 44211
$data[123,4] = 9
$data[123,5] = 100112037
$data[123,6] = 'Cebollín'
$data[123,7] = 'Sin especificar'
$data[123,8] = 'Primera'
$data[123,9] = 110
$data[123,10] = 8000
$data[123,11] = 8000
$data[123,12] = 8000
$data[123,13] = '$/docena de paquetes'
$data[123,14] = 'Provincia de Cautín'
$data[123,15] = 667
$data[123,16] = 12
$data[123,17] = 'Hortaliza'

$data[124,0] = 10
$data[124,1] = 'Vega Modelo de Temuco'
$data[124,2] = 'La Araucanía'
$data[124,3] = 44425
$data[124,4] = 9
$data[124,5] = 100112037
$data[124,6] = 'Cebollín'
$data[124,7] = 'Sin especificar'
$data[124,8] = 'Primera'
$data[124,9] = 30
$data[124,10] = 7000
$data[124,11] = 7000
$data[124,12] = 7000
$data[124,13] = '$/docena de paquetes'
$data[124,14] = 'Provincia de Cautín'
$data[124,15] = 583
$data[124,16] = 12
$data[124,17] = 'Hortaliza'

$data[125,0] = 10
$data[125,1] = 'Vega Modelo de Temuco'
$data[125,2] = 'La Araucanía'
$data[125,3] = 44425
$data[125,4] = 9
$data[125,5] = 100112037
$data[125,6] = 'Cebollín'
$data[125,7] = 'Sin especificar'
$data[125,8] = 'Primera'
$data[125,9] = 20
$data[125,10] = 5000
$data[125,11] = 5000
$data[125,12] = 5000
$data[125,13] = '$/docena de paquetes'
$data[125,14] = 'Región de O''Higgins'
$data[125,15] = 417
$data[125,16] = 12
$data[125,17] = 'Hortaliza'

$ws.Range("A86:R211").Value = $data

# Fix date NumberFormat for the two brand-new rows (210, 211) so column D matches the date style used elsewhere
$dateFormat = $ws.Range("D209").NumberFormat
$ws.Range("D210").NumberFormat = $dateFormat
$ws.Range("D211").NumberFormat = $dateFormat
